$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 10.75879266666667
$ws.Cells.Item(2, 8).Value = 32.276378
$ws.Cells.Item(2, 9).Value = 0.9290725491349732
$ws.Cells.Item(2, 10).Value = 0.9290725491349733
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.971099
$ws.Cells.Item(2, 14).Value = 2.913297
$ws.Cells.Item(2, 15).Value = 0.007882574716876797
$ws.Cells.Item(2, 16).Value = 0.007882574716876797
$ws.Cells.Item(2, 17).Value = 10.44785279980733
$ws.Cells.Item(2, 18).Value = 94.03067519826601
$ws.Cells.Item(2, 19).Value = 0.007323483785955615
$ws.Cells.Item(2, 20).Value = 0.007323483785955616
$ws.Cells.Item(3, 7).Value = 10.75879266666667
$ws.Cells.Item(3, 8).Value = 32.276378
$ws.Cells.Item(3, 9).Value = 0.9290725491349732
$ws.Cells.Item(3, 10).Value = 0.9290725491349733
$ws.Cells.Item(3, 15).Value = 0.6966643430097871
$ws.Cells.Item(3, 16).Value = 0.696664343009787
$ws.Cells.Item(3, 17).Value = 923.3843976203574
$ws.Cells.Item(3, 18).Value = 8310.459578583217
$ws.Cells.Item(3, 19).Value = 0.6472517170515442
$ws.Cells.Item(3, 20).Value = 0.6472517170515442
$ws.Cells.Item(4, 7).Value = 10.75879266666667
$ws.Cells.Item(4, 8).Value = 32.276378
$ws.Cells.Item(4, 9).Value = 0.9290725491349732
$ws.Cells.Item(4, 10).Value = 0.9290725491349733
$ws.Cells.Item(4, 13).Value = 36.24916566666667
$ws.Cells.Item(4, 14).Value = 108.747497
$ws.Cells.Item(4, 15).Value = 0.294240604502677
$ws.Cells.Item(4, 16).Value = 0.294240604502677
$ws.Cells.Item(4, 17).Value = 389.9972577473185
$ws.Cells.Item(4, 18).Value = 3509.975319725866
$ws.Cells.Item(4, 19).Value = 0.2733708684843176
$ws.Cells.Item(4, 20).Value = 0.2733708684843176
$ws.Cells.Item(5, 7).Value = 10.75879266666667
$ws.Cells.Item(5, 8).Value = 32.276378
$ws.Cells.Item(5, 9).Value = 0.9290725491349732
$ws.Cells.Item(5, 10).Value = 0.9290725491349733
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.149372
$ws.Cells.Item(5, 14).Value = 0.448116
$ws.Cells.Item(5, 15).Value = 0.001212477770659141
$ws.Cells.Item(5, 16).Value = 0.001212477770659141
$ws.Cells.Item(5, 17).Value = 1.607062378205333
$ws.Cells.Item(5, 18).Value = 14.463561403848
$ws.Cells.Item(5, 19).Value = 0.001126479813155777
$ws.Cells.Item(5, 20).Value = 0.001126479813155777
$ws.Cells.Item(6, 8).Value = 0.627738
$ws.Cells.Item(6, 9).Value = 0.01806938014695731
$ws.Cells.Item(6, 10).Value = 0.01806938014695731
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.971099
$ws.Cells.Item(6, 14).Value = 2.913297
$ws.Cells.Item(6, 15).Value = 0.007882574716876797
$ws.Cells.Item(6, 16).Value = 0.007882574716876797
$ws.Cells.Item(6, 17).Value = 0.203198581354
$ws.Cells.Item(6, 18).Value = 1.828787232186
$ws.Cells.Item(6, 19).Value = 0.0001424332390960413
$ws.Cells.Item(6, 20).Value = 0.0001424332390960413
$ws.Cells.Item(7, 8).Value = 0.627738
$ws.Cells.Item(7, 9).Value = 0.01806938014695731
$ws.Cells.Item(7, 10).Value = 0.01806938014695731
$ws.Cells.Item(7, 15).Value = 0.6966643430097871
$ws.Cells.Item(7, 16).Value = 0.696664343009787
$ws.Cells.Item(7, 19).Value = 0.01258829284867411
$ws.Cells.Item(7, 20).Value = 0.0125882928486741
$ws.Cells.Item(8, 8).Value = 0.627738
$ws.Cells.Item(8, 9).Value = 0.01806938014695731
$ws.Cells.Item(8, 10).Value = 0.01806938014695731
$ws.Cells.Item(8, 13).Value = 36.24916566666667
$ws.Cells.Item(8, 14).Value = 108.747497
$ws.Cells.Item(8, 15).Value = 0.294240604502677
$ws.Cells.Item(8, 16).Value = 0.294240604502677
$ws.Cells.Item(8, 17).Value = 7.584992919087335
$ws.Cells.Item(8, 18).Value = 68.26493627178601
$ws.Cells.Item(8, 19).Value = 0.005316745337429391
$ws.Cells.Item(8, 20).Value = 0.00531674533742939
$ws.Cells.Item(9, 8).Value = 0.627738
$ws.Cells.Item(9, 9).Value = 0.01806938014695731
$ws.Cells.Item(9, 10).Value = 0.01806938014695731
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.149372
$ws.Cells.Item(9, 14).Value = 0.448116
$ws.Cells.Item(9, 15).Value = 0.001212477770659141
$ws.Cells.Item(9, 16).Value = 0.001212477770659141
$ws.Cells.Item(9, 17).Value = 0.031255493512
$ws.Cells.Item(9, 18).Value = 0.281299441608
$ws.Cells.Item(9, 19).Value = 0.00002190872175777534
$ws.Cells.Item(9, 20).Value = 0.00002190872175777534
$ws.Cells.Item(10, 7).Value = 0.5786906666666666
$ws.Cells.Item(10, 8).Value = 1.736072
$ws.Cells.Item(10, 9).Value = 0.04997267160899686
$ws.Cells.Item(10, 10).Value = 0.04997267160899686
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.971099
$ws.Cells.Item(10, 14).Value = 2.913297
$ws.Cells.Item(10, 15).Value = 0.007882574716876797
$ws.Cells.Item(10, 16).Value = 0.007882574716876797
$ws.Cells.Item(10, 17).Value = 0.5619659277093333
$ws.Cells.Item(10, 18).Value = 5.057693349384
$ws.Cells.Item(10, 19).Value = 0.0003939133177598656
$ws.Cells.Item(10, 20).Value = 0.0003939133177598656
$ws.Cells.Item(11, 7).Value = 0.5786906666666666
$ws.Cells.Item(11, 8).Value = 1.736072
$ws.Cells.Item(11, 9).Value = 0.04997267160899686
$ws.Cells.Item(11, 10).Value = 0.04997267160899686
$ws.Cells.Item(11, 15).Value = 0.6966643430097871
$ws.Cells.Item(11, 16).Value = 0.696664343009787
$ws.Cells.Item(11, 17).Value = 49.66671904590932
$ws.Cells.Item(11, 18).Value = 447.000471413184
$ws.Cells.Item(11, 19).Value = 0.03481417843492564
$ws.Cells.Item(11, 20).Value = 0.03481417843492563
$ws.Cells.Item(12, 7).Value = 0.5786906666666666
$ws.Cells.Item(12, 8).Value = 1.736072
$ws.Cells.Item(12, 9).Value = 0.04997267160899686
$ws.Cells.Item(12, 10).Value = 0.04997267160899686
$ws.Cells.Item(12, 13).Value = 36.24916566666667
$ws.Cells.Item(12, 14).Value = 108.747497
$ws.Cells.Item(12, 15).Value = 0.294240604502677
$ws.Cells.Item(12, 16).Value = 0.294240604502677
$ws.Cells.Item(12, 17).Value = 20.97705384575378
$ws.Cells.Item(12, 18).Value = 188.793484611784
$ws.Cells.Item(12, 19).Value = 0.014703989102845
$ws.Cells.Item(12, 20).Value = 0.014703989102845
$ws.Cells.Item(13, 7).Value = 0.5786906666666666
$ws.Cells.Item(13, 8).Value = 1.736072
$ws.Cells.Item(13, 9).Value = 0.04997267160899686
$ws.Cells.Item(13, 10).Value = 0.04997267160899686
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.149372
$ws.Cells.Item(13, 14).Value = 0.448116
$ws.Cells.Item(13, 15).Value = 0.001212477770659141
$ws.Cells.Item(13, 16).Value = 0.001212477770659141
$ws.Cells.Item(13, 17).Value = 0.08644018226133331
$ws.Cells.Item(13, 18).Value = 0.7779616403519999
$ws.Cells.Item(13, 19).Value = 0.00006059075346635785
$ws.Cells.Item(13, 20).Value = 0.00006059075346635785
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.03341333333333333
$ws.Cells.Item(14, 8).Value = 0.10024
$ws.Cells.Item(14, 9).Value = 0.002885399109072576
$ws.Cells.Item(14, 10).Value = 0.002885399109072576
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.971099
$ws.Cells.Item(14, 14).Value = 2.913297
$ws.Cells.Item(14, 15).Value = 0.007882574716876797
$ws.Cells.Item(14, 16).Value = 0.007882574716876797
$ws.Cells.Item(14, 17).Value = 0.03244765458666667
$ws.Cells.Item(14, 18).Value = 0.29202889128
$ws.Cells.Item(14, 19).Value = 0.00002274437406527433
$ws.Cells.Item(14, 20).Value = 0.00002274437406527433
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.03341333333333333
$ws.Cells.Item(15, 8).Value = 0.10024
$ws.Cells.Item(15, 9).Value = 0.002885399109072576
$ws.Cells.Item(15, 10).Value = 0.002885399109072576
$ws.Cells.Item(15, 15).Value = 0.6966643430097871
$ws.Cells.Item(15, 16).Value = 0.696664343009787
$ws.Cells.Item(15, 17).Value = 2.867733548586667
$ws.Cells.Item(15, 18).Value = 25.80960193728
$ws.Cells.Item(15, 19).Value = 0.002010154674643072
$ws.Cells.Item(15, 20).Value = 0.002010154674643071
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.03341333333333333
$ws.Cells.Item(16, 8).Value = 0.10024
$ws.Cells.Item(16, 9).Value = 0.002885399109072576
$ws.Cells.Item(16, 10).Value = 0.002885399109072576
$ws.Cells.Item(16, 13).Value = 36.24916566666667
$ws.Cells.Item(16, 14).Value = 108.747497
$ws.Cells.Item(16, 15).Value = 0.294240604502677
$ws.Cells.Item(16, 16).Value = 0.294240604502677
$ws.Cells.Item(16, 17).Value = 1.211205455475556
$ws.Cells.Item(16, 18).Value = 10.90084909928
$ws.Cells.Item(16, 19).Value = 0.0008490015780850006
$ws.Cells.Item(16, 20).Value = 0.0008490015780850004
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.03341333333333333
$ws.Cells.Item(17, 8).Value = 0.10024
$ws.Cells.Item(17, 9).Value = 0.002885399109072576
$ws.Cells.Item(17, 10).Value = 0.002885399109072576
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.149372
$ws.Cells.Item(17, 14).Value = 0.448116
$ws.Cells.Item(17, 15).Value = 0.001212477770659141
$ws.Cells.Item(17, 16).Value = 0.001212477770659141
$ws.Cells.Item(17, 17).Value = 0.004991016426666666
$ws.Cells.Item(17, 18).Value = 0.04491914783999999
$ws.Cells.Item(17, 19).Value = 0.000003498482279230188
$ws.Cells.Item(17, 20).Value = 0.000003498482279230188
